$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D8").Value = "04-reading.html"
$ws.Range("D9").Value = "05-reading.html"
$ws.Range("D11").Value = "06-reading.html"
$ws.Range("D14").Value = "08-reading.html"
$ws.Range("D14").Select()
